# Auto-generated Excel COM-interop edit script
# Applies updated crypto price/volume data per Thu Oct 19 03:35:00 UTC 2023 GitHub Actions run

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '28.288.89'
$ws.Range("E2").Value = '  -0.85%  '

$ws.Range("D3").Value = '1.551.93'
$ws.Range("E3").Value = '  -1.02%  '

$ws.Range("E4").Value = '  -0.06%  '

$ws.Range("E5").Value = '  -1.60%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.483'
$ws.Range("E6").Value = '  -1.77%  '

$ws.Range("E7").Value = '  -0.08%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '23.43'
$ws.Range("E8").Value = '  -2.70%  '

$ws.Range("E9").Value = '  -2.18%  '

$ws.Range("E10").Value = '  -1.28%  '

$ws.Range("E11").Value = '  +0.05%  '

$ws.Range("D12").Value = '1.773.82'
$ws.Range("E12").Value = '  -1.00%  '

$ws.Range("D13").Value = '1.551.63'
$ws.Range("E13").Value = '  -1.02%  '

$ws.Range("D14").Value = '28.296.05'
$ws.Range("E14").Value = '  -0.76%  '

$ws.Range("E16").Value = '  -2.35%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '60.45'
$ws.Range("E17").Value = '  -2.80%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '226.47'
$ws.Range("E18").Value = '  -1.61%  '

$ws.Range("E19").Value = '  -0.66%  '

$ws.Range("E21").Value = '  -0.06%  '

$ws.Range("E22").Value = '  +1.17%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '8.80'
$ws.Range("E23").Value = '  -3.36%  '

$ws.Range("E24").Value = '  -5.35%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '147.91'
$ws.Range("E25").Value = '  -2.25%  '

$ws.Range("E26").Value = '  -1.69%  '

$ws.Range("E27").Value = '  -0.31%  '

$ws.Range("E28").Value = '  -0.04%  '

$ws.Range("E29").Value = '  -3.07%  '

$ws.Range("E30").Value = '  -3.66%  '

$ws.Range("E31").Value = '  -4.13%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.17'
$ws.Range("E32").Value = '  -0.83%  '

$ws.Range("E33").Value = '  -1.04%  '

$ws.Range("D34").Value = '1.385.70'
$ws.Range("E34").Value = '  -0.48%  '

$ws.Range("E35").Value = '  +0.36%  '

$ws.Range("E36").Value = '  -2.65%  '

$ws.Range("E37").Value = '  -1.51%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.59'
$ws.Range("E38").Value = '  -1.23%  '

$ws.Range("E39").Value = '  -2.42%  '

$ws.Range("E40").Value = '  +1.82%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.511'
$ws.Range("E41").Value = '  -2.09%  '

$ws.Range("E42").Value = '  -0.05%  '

$ws.Range("E43").Value = '  -1.46%  '

$ws.Range("E44").Value = '  +0.75%  '

$ws.Range("E45").Value = '  -1.11%  '

$ws.Range("E46").Value = '  -1.56%  '

$ws.Range("E47").Value = '  -1.06%  '

$ws.Range("E48").Value = '  -6.64%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '85.40'
$ws.Range("E49").Value = '  -1.08%  '

# Row 50 and 51 swap: BabyDogeCoin/BitcoinSV order reversed with updated values
$ws.Range("B50").Value = 'BitcoinSV'
$ws.Range("C50").Value = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '41.97'
$ws.Range("E50").Value = '  +6.04%  '

$ws.Range("B51").Value = 'BabyDogeCoin'
$ws.Range("C51").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D51").Value = '0.0₆0103'
$ws.Range("E51").Value = '  +0.07%  '
